# ============================================================================
# docs: add medium-scale M3 demo results (runs 8-9) to DOCX report
# ============================================================================

$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# 1) Key Metrics table (first table): "25 p/s RabbitMQ," / "27 p/s Redis"
#    -> "256 p/s RabbitMQ," / "243 p/s Redis (10K)"
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("25 p/s RabbitMQ,", $true, $false, $false, $false, $false, `
    $true, 1, $false, "256 p/s RabbitMQ,", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("256 p/s RabbitMQ,", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.MoveEnd(1, 12)
$rng.Text = "243 p/s Redis (10K)"

# ----------------------------------------------------------------------------
# 2) M3 comparison table "Parameters Processed" row, demo column:
#    "1,000 (demo)" -> "10,000"
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("1,000 (demo)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "10,000", 2) | Out-Null

# ----------------------------------------------------------------------------
# 3) M3 comparison table "Execution Time" row, demo column:
#    "39s (49s with fault" / "demo)" -> "~39s (RMQ) / ~41s" / "(Redis)"
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("39s (49s with fault", $true, $false, $false, $false, $false, `
    $true, 1, $false, "~39s (RMQ) / ~41s", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("~39s (RMQ) / ~41s", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.MoveEnd(1, 6)
$rng.Text = "(Redis)"

# ----------------------------------------------------------------------------
# 4) M3 comparison table "Throughput" row, demo column:
#    "25 params/sec (demo" / "scale)" -> "256 p/s (RMQ)" / "243 p/s (Redis)"
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("25 params/sec (demo", $true, $false, $false, $false, $false, `
    $true, 1, $false, "256 p/s (RMQ)", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("256 p/s (RMQ)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.MoveEnd(1, 6)
$rng.Text = "243 p/s (Redis)"

# ----------------------------------------------------------------------------
# 5) "100% success rate across all 7 verified demo runs" -> "...all 9..."
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("100% success rate across all 7 verified demo runs", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "100% success rate across all 9 verified demo runs", 2) | Out-Null

# ----------------------------------------------------------------------------
# 6) Insert two new Key Takeaways bullets right after the "Throughput drops
#    ~20-25%..." bullet and before the "RabbitMQ adds retry queues..." bullet.
# ----------------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Throughput drops*") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($targetIndex + 1)
$newPara1.Range.Text = "•  At equal scale (10K params, medium preset), RabbitMQ (256 p/s) slightly outperforms Redis Streams (243 p/s) — both backends are comparable."

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($targetIndex + 2)
$newPara2.Range.Text = "•  Runs 4-7 used the 'small' preset (1K params, 100ms simulated work per chunk), explaining the lower throughput vs. runs 1-3 and 8-9 (10ms per chunk)."

# ----------------------------------------------------------------------------
# 7) Final bullet: add "--scale small|medium" to the demo command.
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("Both backends can be demoed with a single command: ./scripts/run-demo.sh --backend redis|rabbitmq", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Both backends can be demoed with a single command: ./scripts/run-demo.sh --backend redis|rabbitmq --scale small|medium", 2) | Out-Null

# ----------------------------------------------------------------------------
# 8) Append two new rows (runs 8 and 9) to the M3 demo-runs table (last table).
# ----------------------------------------------------------------------------
$t = $d.Tables.Item($d.Tables.Count)
$lf = [char]11

$row8 = $t.Rows.Add()
$row8.Cells.Item(1).Range.Text = "8"
$row8.Cells.Item(2).Range.Text = "M3"
$row8.Cells.Item(3).Range.Text = "Redis" + $lf + "Streams"
$row8.Cells.Item(4).Range.Text = "No"
$row8.Cells.Item(5).Range.Text = "3"
$row8.Cells.Item(6).Range.Text = "10,000"
$row8.Cells.Item(7).Range.Text = "100/100"
$row8.Cells.Item(8).Range.Text = "41s"
$row8.Cells.Item(9).Range.Text = "243 p/s"
$row8.Cells.Item(10).Range.Text = "ZERO"

$row9 = $t.Rows.Add()
$row9.Cells.Item(1).Range.Text = "9"
$row9.Cells.Item(2).Range.Text = "M3"
$row9.Cells.Item(3).Range.Text = "RabbitMQ"
$row9.Cells.Item(4).Range.Text = "No"
$row9.Cells.Item(5).Range.Text = "3"
$row9.Cells.Item(6).Range.Text = "10,000"
$row9.Cells.Item(7).Range.Text = "100/100"
$row9.Cells.Item(8).Range.Text = "38s"
$row9.Cells.Item(9).Range.Text = "256 p/s"
$row9.Cells.Item(10).Range.Text = "ZERO"

Write-Output "edit complete"
